$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.918.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = "'1.898.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'0.7942"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.95%  '
$ws.Range("D6").Value = "'244.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = "'0.3172"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.58%  '
$ws.Range("D9").Value = "'25.57"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.68%  '
$ws.Range("D10").Value = "'0.07192"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.10%  '
$ws.Range("D11").Value = "'0.08114"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'5.654"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.80%  '
$ws.Range("D13").Value = "'0.7712"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.37%  '
$ws.Range("D14").Value = "'1.890.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("D15").Value = "'92.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").Value = "'6.179"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.44%  '
$ws.Range("D17").Value = "'29.908.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").Value = "'13.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("D19").Value = "'245.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("D20").Value = "'0.000007776"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").Value = "'8.258"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +18.50%  '
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = "'2.143.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").Value = "'0.1675"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.72%  '
$ws.Range("D26").Value = "'9.523"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.02%  '
$ws.Range("D27").Value = "'164.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.92%  '
$ws.Range("D28").Value = "'18.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.42%  '
$ws.Range("D29").Value = "'2.078"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("D30").Value = "'1.405"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.19%  '
$ws.Range("D31").Value = "'1.550"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.07%  '
$ws.Range("D32").Value = "'4.507"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.21%  '
$ws.Range("D33").Value = "'0.05646"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.32%  '
$ws.Range("D34").Value = "'4.097"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.56%  '
$ws.Range("D35").Value = "'1.287"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.63%  '
$ws.Range("D36").Value = "'0.7459"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.32%  '
$ws.Range("D37").Value = "'0.9984"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("D38").Value = "'2.631"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.41%  '
$ws.Range("D39").Value = "'0.01939"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("E40").Value = '  +0.36%  '
$ws.Range("D41").Value = "'1.166.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +15.96%  '
$ws.Range("D42").Value = "'75.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.63%  '
$ws.Range("D43").Value = "'0.4441"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.27%  '
$ws.Range("E44").Value = '  +1.53%  '
$ws.Range("D45").Value = "'0.8542"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("D46").Value = "'104.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.95%  '
$ws.Range("D47").Value = "'1.000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("D48").Value = "'1.892"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("D49").Value = "'10.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.16%  '
$ws.Range("D50").Value = "'7.515"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("D51").Value = "'3.016"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +10.61%  '
